# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# 1. About sheet: bump the "Source:" date (C1) from 2024-03-15 to 2024-03-28.
# 2. RAF-capacity sheet: raise the capacity-credit multiplier for the two
#    hydrogen technologies (hydrogen combustion turbine / hydrogen combined
#    cycle, rows 24-25) from 0.3 to 1.
# 3. View-state touch-ups: RAF-capacity becomes the active/selected sheet
#    (instead of RAF-generation), its zoom is set to 80%, and the last
#    selected cell there is B25. Its first column is also widened to fit
#    the longer hydrogen technology labels.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsGen   = $wb.Worksheets.Item("RAF-generation")
$wsDAT   = $wb.Worksheets.Item("RAF-demand-altering-techs")
$wsCap   = $wb.Worksheets.Item("RAF-capacity")

# --- 1. Update the source date on the About sheet ------------------------
$wsAbout.Range("C1").Value = 45379

# --- 2. Update hydrogen capacity-credit multipliers -----------------------
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1

# --- 3. Widen column A on RAF-capacity so the hydrogen labels fit ---------
$wsCap.Columns.Item(1).ColumnWidth = 28.1667

# --- 4. Make RAF-capacity the active sheet / tab, at 80% zoom, with the
#        last selection on B25 (mirrors how RAF-generation previously had
#        tabSelected + activeTab pointing at it) ---------------------------
$wsCap.Activate()
$excel.ActiveWindow.Zoom = 80
$wsCap.Range("B25").Select()
